$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.472.74"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.898.44"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4914"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2933"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06707"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "1.915.58"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.160"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6658"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "30.446.99"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007844"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "2.138.52"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.359"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.115"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.491"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.944"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.469"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.344"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.049"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05194"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7420"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.101"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01811"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.673"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9229"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.044"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4388"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.914"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9944"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +20.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1372"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.608"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.051"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05833"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3928"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
